# Applies targeted numeric corrections to several Leve profit rows
# across multiple job sheets (ALC, ARM, BSM, CRP, GSM, LTW, WVR),
# refreshing currentAveragePrice / LevePrice / LeveProfit columns.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H55").Value = 122
$ws.Range("J55").Value = 122
$ws.Range("L55").Value = 122
$ws.Range("N55").Value = -550
$ws.Range("H74").Value = 6750
$ws.Range("I74").Value = 6750
$ws.Range("K74").Value = 6750
$ws.Range("M74").Value = -5814
$ws.Range("H77").Value = 6750
$ws.Range("I77").Value = 6750
$ws.Range("K77").Value = 33750
$ws.Range("M77").Value = -29070
$ws.Range("H106").Value = 5100
$ws.Range("I106").Value = 5100
$ws.Range("K106").Value = 5100
$ws.Range("M106").Value = -4469
$ws.Range("H111").Value = 716
$ws.Range("I111").Value = 500
$ws.Range("J111").Value = 932
$ws.Range("K111").Value = 1500
$ws.Range("L111").Value = 2796
$ws.Range("M111").Value = 1567
$ws.Range("N111").Value = -8930

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 0
$ws.Range("I45").Value = 0
$ws.Range("K45").Value = 0
$ws.Range("M45").ClearContents()
$ws.Range("H61").Value = 2950
$ws.Range("I61").Value = 2331.6924
$ws.Range("K61").Value = 2331.6924
$ws.Range("M61").Value = -2119.6924
$ws.Range("H132").Value = 3704.5715
$ws.Range("I132").Value = 2655.3333
$ws.Range("K132").Value = 7965.999899999999
$ws.Range("M132").Value = -5435.999899999999
$ws.Range("H136").Value = 2950
$ws.Range("I136").Value = 2331.6924
$ws.Range("K136").Value = 6995.0772
$ws.Range("M136").Value = -4445.0772

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 8166.6665
$ws.Range("I20").Value = 10000
$ws.Range("K20").Value = 10000
$ws.Range("M20").Value = -9753
$ws.Range("H22").Value = 308.77777
$ws.Range("I22").Value = 308.77777
$ws.Range("K22").Value = 308.77777
$ws.Range("M22").Value = -135.77777

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 9999.5
$ws.Range("I16").Value = 9999.5
$ws.Range("K16").Value = 9999.5
$ws.Range("M16").Value = -9712.5
$ws.Range("H22").Value = 918.2
$ws.Range("I22").Value = 922.75
$ws.Range("J22").Value = 900
$ws.Range("K22").Value = 922.75
$ws.Range("L22").Value = 900
$ws.Range("M22").Value = -572.75
$ws.Range("N22").Value = -1600
$ws.Range("H58").Value = 6152.4287
$ws.Range("I58").Value = 4385.4165
$ws.Range("K58").Value = 4385.4165
$ws.Range("M58").Value = -4182.4165
$ws.Range("H113").Value = 9999.5
$ws.Range("I113").Value = 9999.5
$ws.Range("K113").Value = 9999.5
$ws.Range("M113").Value = -7829.5
$ws.Range("H122").Value = 4055.9412
$ws.Range("I122").Value = 4055.9412
$ws.Range("J122").Value = 0
$ws.Range("K122").Value = 12167.8236
$ws.Range("L122").Value = 0
$ws.Range("M122").Value = -9717.8236
$ws.Range("N122").ClearContents()
$ws.Range("H132").Value = 827.1875
$ws.Range("I132").Value = 827.1875
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 2481.5625
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = 48.4375
$ws.Range("N132").ClearContents()
$ws.Range("H136").Value = 6152.4287
$ws.Range("I136").Value = 4385.4165
$ws.Range("K136").Value = 13156.2495
$ws.Range("M136").Value = -10606.2495

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 6099.5386
$ws.Range("I70").Value = 5660.125
$ws.Range("J70").Value = 6802.6
$ws.Range("K70").Value = 5660.125
$ws.Range("L70").Value = 6802.6
$ws.Range("M70").Value = -5390.125
$ws.Range("N70").Value = -7342.6
$ws.Range("H73").Value = 6099.5386
$ws.Range("I73").Value = 5660.125
$ws.Range("J73").Value = 6802.6
$ws.Range("K73").Value = 5660.125
$ws.Range("L73").Value = 6802.6
$ws.Range("M73").Value = -4724.125
$ws.Range("N73").Value = -8674.6
$ws.Range("I99").Value = 24471
$ws.Range("J99").Value = 35000
$ws.Range("K99").Value = 24471
$ws.Range("L99").Value = 35000
$ws.Range("M99").Value = -22225
$ws.Range("N99").Value = -39492
$ws.Range("H113").Value = 1203.2
$ws.Range("I113").Value = 1166.5
$ws.Range("K113").Value = 1166.5
$ws.Range("M113").Value = 1003.5
$ws.Range("H122").Value = 46609.266
$ws.Range("J122").Value = 37834.332
$ws.Range("L122").Value = 113502.996
$ws.Range("N122").Value = -118402.996

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 1200
$ws.Range("I22").Value = 1200
$ws.Range("K22").Value = 1200
$ws.Range("M22").Value = -905
$ws.Range("H27").Value = 1200
$ws.Range("I27").Value = 1200
$ws.Range("K27").Value = 1200
$ws.Range("M27").Value = -1093
$ws.Range("H61").Value = 4650
$ws.Range("I61").Value = 4650
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 4650
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -4448
$ws.Range("N61").ClearContents()
$ws.Range("H100").Value = 6185.875
$ws.Range("I100").Value = 6185.875
$ws.Range("K100").Value = 6185.875
$ws.Range("M100").Value = -5644.875
$ws.Range("H113").Value = 4650
$ws.Range("I113").Value = 4650
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 4650
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -2480
$ws.Range("N113").ClearContents()
$ws.Range("H132").Value = 18133.572
$ws.Range("I132").Value = 18987
$ws.Range("J132").Value = 16000
$ws.Range("K132").Value = 56961
$ws.Range("L132").Value = 48000
$ws.Range("M132").Value = -54431
$ws.Range("N132").Value = -53060
$ws.Range("H136").Value = 6345.222
$ws.Range("I136").Value = 3150.75
$ws.Range("K136").Value = 9452.25
$ws.Range("M136").Value = -6902.25

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H81").Value = 5162.5
$ws.Range("I81").Value = 925
$ws.Range("J81").Value = 9400
$ws.Range("K81").Value = 1850
$ws.Range("L81").Value = 18800
$ws.Range("M81").Value = -789
$ws.Range("N81").Value = -20922
$ws.Range("H84").Value = 5162.5
$ws.Range("I84").Value = 925
$ws.Range("J84").Value = 9400
$ws.Range("K84").Value = 9250
$ws.Range("L84").Value = 94000
$ws.Range("M84").Value = -3946
$ws.Range("N84").Value = -104608
$ws.Range("H107").Value = 438.6
$ws.Range("I107").Value = 431
$ws.Range("J107").Value = 450
$ws.Range("K107").Value = 1293
$ws.Range("L107").Value = 1350
$ws.Range("M107").Value = 627
$ws.Range("N107").Value = -5190
$ws.Range("H113").Value = 1536.8096
$ws.Range("I113").Value = 839.375
$ws.Range("J113").Value = 3768.6
$ws.Range("K113").Value = 2518.125
$ws.Range("L113").Value = 11305.8
$ws.Range("M113").Value = -348.125
$ws.Range("N113").Value = -15645.8
$ws.Range("H136").Value = 3069.75
$ws.Range("I136").Value = 3069.75
$ws.Range("K136").Value = 9209.25
$ws.Range("M136").Value = -6659.25
